# Inflation driver input workbook: add a "Sheet1" validation-list sheet and
# rework the Tabelle1 parameter sheet (new row labels + a "Currency"/"ZC
# Inflation" driver column B) with a dropdown validation on B3:AD3.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Re-label column A of Tabelle1 for the new risk-driver fields.
# ---------------------------------------------------------------------
$ws1.Range("A7").Value = "Volatility Inflation Index"
$ws1.Range("A8").Value = "Mean reversion Real IR"
$ws1.Range("A6").Value = "Volatility Real IR"
$ws1.Range("A3").Value = "Currency"
$ws1.Range("A4").Value = "ZC Inflation"

# ---------------------------------------------------------------------
# 2) Add the new "Sheet1" tab right after "Tabelle1" and fill its B1:B27
#    currency-code list (this is the source range for the dropdown).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

$ws2.Range("B1").Value = "domestic"
$ws2.Range("B2").Value = "foreign1"
$ws2.Range("B3").Value = "foreign2"
$ws2.Range("B4").Value = "foreign3"
$ws2.Range("B5").Value = "foreign4"
$ws2.Range("B6").Value = "foreign5"
$ws2.Range("B7").Value = "foreign6"
$ws2.Range("B8").Value = "foreign7"
$ws2.Range("B9").Value = "foreign8"
$ws2.Range("B10").Value = "foreign9"
$ws2.Range("B11").Value = "foreign10"
$ws2.Range("B12").Value = "foreign11"
$ws2.Range("B13").Value = "foreign12"
$ws2.Range("B14").Value = "foreign13"
$ws2.Range("B15").Value = "foreign14"
$ws2.Range("B16").Value = "foreign15"
$ws2.Range("B17").Value = "foreign16"
$ws2.Range("B18").Value = "foreign17"
$ws2.Range("B19").Value = "foreign18"
$ws2.Range("B20").Value = "foreign19"
$ws2.Range("B21").Value = "foreign20"
$ws2.Range("B22").Value = "foreign21"
$ws2.Range("B23").Value = "foreign22"
$ws2.Range("B24").Value = "foreign23"
$ws2.Range("B25").Value = "foreign24"
$ws2.Range("B26").Value = "foreign25"
$ws2.Range("B27").Value = "foreign26"

# ---------------------------------------------------------------------
# 3) Fill in the sample driver values in column B.
# ---------------------------------------------------------------------
$ws1.Range("B2").Value = "hicp"
$ws1.Range("B3").Value = "foreign1"
$ws1.Range("B4").Value = "EUR Inflation 31122019"
$ws1.Range("B5").Value = 105
$ws1.Range("B6").Value = "EUR Real Vol"
$ws1.Range("B7").Value = "EUR Inflation Vol"
$ws1.Range("B8").Value = 0.03
$ws1.Range("B8").NumberFormat = "0%"

# ---------------------------------------------------------------------
# 4) Dropdown (list) validation on B3:AD3, sourced from Sheet1!$B$1:$B$27.
# ---------------------------------------------------------------------
$dvRange = $ws1.Range("B3:AD3")
$dvRange.Validation.Delete()
$dvRange.Validation.Add(3, 1, 1, '=Sheet1!$B$1:$B$27')
$dvRange.Validation.IgnoreBlank = $true
$dvRange.Validation.InCellDropdown = $true
$dvRange.Validation.ShowInput = $true
$dvRange.Validation.ShowError = $true

# ---------------------------------------------------------------------
# 5) Restore the on-screen selections seen in the authored workbook.
# ---------------------------------------------------------------------
[void]$ws2.Range("B2:B27").Select()
[void]$ws1.Range("B15").Select()
[void]$ws1.Activate()
